# Update "想去人数" (want-to-go count, column F) figures to the values
# captured in the newer data pull (gh-pages output generated at 456a3b4).
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 3 = 本地生活 (Local life)
# Sheet 4 = 全部类型 (All types - aggregates rows from every other sheet,
#           so several of the same events/values reappear here)

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value  = 8037   # 杭州·2024首届COMIC GALAXY次元盛典
$ws.Range("F4").Value  = 1894   # 杭州·浮游猫动漫嘉年华
$ws.Range("F5").Value  = 6474   # 杭州·理想乡动漫展-同人创作者大会
$ws.Range("F7").Value  = 2024   # 杭州·Eternal时光国乙only展（日+夜场）
$ws.Range("F8").Value  = 557    # 杭州·第五人格同人only
$ws.Range("F15").Value = 8383   # 杭州·第二届次元格子动漫展-进入格子空间，探索次元世界！
$ws.Range("F21").Value = 857    # 杭州·萌忧 原崩铁同人only
$ws.Range("F29").Value = 1992   # 杭州·首届CCPC动漫嘉年华

# --- 本地生活 (sheet 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value  = 2312   # 杭州·盗墓笔记官方授权「四季同书」主题店

# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value  = 2312   # 杭州·盗墓笔记官方授权「四季同书」主题店
$ws.Range("F6").Value  = 8037   # 杭州·2024首届COMIC GALAXY次元盛典
$ws.Range("F9").Value  = 1894   # 杭州·浮游猫动漫嘉年华
$ws.Range("F10").Value = 6474   # 杭州·理想乡动漫展-同人创作者大会
$ws.Range("F11").Value = 2024   # 杭州·Eternal时光国乙only展（日+夜场）
$ws.Range("F13").Value = 557    # 杭州·第五人格同人only
$ws.Range("F23").Value = 8383   # 杭州·第二届次元格子动漫展-进入格子空间，探索次元世界！
$ws.Range("F29").Value = 857    # 杭州·萌忧 原崩铁同人only
$ws.Range("F35").Value = 1992   # 杭州·首届CCPC动漫嘉年华
